$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): values only; style is already "1" (header/gray) for
# every existing cell A1:H1, so we just need to re-target text and add two
# more header cells (I1, J1) with the same header style.
$ws.Range("C1").Value = "Meal type"
$ws.Range("D1").Value = "Description"
$ws.Range("E1").Value = "Carbs"
$ws.Range("F1").Value = "Glycemic index"
$ws.Range("G1").Value = "Insulin (units)"
$ws.Range("H1").Value = "High correction (units)"
$ws.Range("I1").Value = "Sports correction (units)"
$ws.Range("J1").Value = "Total insulin (units)"

# New header cells need the header format copied onto them.
$ws.Range("A1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Row 2 (Breakfast / Chicken) ---
$ws.Range("C2").Value = "Breakfast"
$ws.Range("D2").Value = "Chicken"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 30
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 2
$ws.Range("J2").Value = 7

# --- Row 3 (Snack / Banana) ---
$ws.Range("C3").Value = "Snack"
$ws.Range("D3").Value = "Banana"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 40
$ws.Range("G3").Value = 4
$ws.Range("H3").ClearContents()
$ws.Range("J3").Value = 4

# --- Row 4 (Dinner / Cheese) ---
$ws.Range("C4").Value = "Dinner"
$ws.Range("D4").Value = "Cheese"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 50
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 1
$ws.Range("J4").Value = 3

# Cells that are brand-new (didn't exist in the original 8-column layout)
# need the plain data-cell format ("s=4") copied onto them; D2/D3's existing
# style already equals that, so it makes a convenient template.
$ws.Range("D2").Copy()
$ws.Range("G2,F3,G3,J2,J3,G4,J4").PasteSpecial(-4122)

# D4 is newly typed text in a cell that used to hold a number (old Carbs);
# the new value should carry NO explicit style, matching a brand-new cell.
$ws.Range("D4").ClearFormats()
$ws.Range("D4").Value = "Cheese"

# --- Column widths shift from F/G/H to H/I/J ---
$ws.Columns.Item("H").ColumnWidth = 19.166666666666668
$ws.Columns.Item("I").ColumnWidth = 21.498697916666668
$ws.Columns.Item("J").ColumnWidth = 19.498697916666668
